# Update April 2024 library statistics (Circulation, ILL Loans, ILL Borrows)
# for each library row in the "Next Statistics" executive board workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Circulation (B), ILL Loans (C), ILL Borrows (D)
$data = @(
    @(3, 21764, 3651, 3686),
    @(4, 12808, 1412, 1583),
    @(5, 36302, 3649, 3509),
    @(6, 610, 379, 80),
    @(7, 24526, 4268, 3236),
    @(8, 3039, 517, 687),
    @(9, 2648, 630, 388),
    @(10, 1210, 211, 113),
    @(11, 330, 172, 2),
    @(12, 0, 0, 0),
    @(13, 494, 105, 175),
    @(14, 1571, 611, 676),
    @(15, 2527, 970, 477),
    @(16, 1788, 888, 210),
    @(17, 1357, 409, 197),
    @(18, 8704, 1301, 1686),
    @(19, 625, 299, 178),
    @(20, 8947, 1465, 1601),
    @(21, 145, 257, 9),
    @(22, 8199, 1101, 1393),
    @(23, 536, 398, 98),
    @(24, 9650, 1066, 2023),
    @(25, 41935, 3668, 5218),
    @(26, 3090, 972, 459),
    @(27, 0, 0, 0),
    @(28, 2462, 662, 604),
    @(29, 937, 318, 224),
    @(30, 7067, 1289, 1644),
    @(31, 271, 62, 159),
    @(32, 1464, 958, 184),
    @(33, 6769, 1757, 1576),
    @(34, 4799, 1792, 972),
    @(35, 2527, 301, 499),
    @(36, 28336, 3338, 2857),
    @(37, 4021, 1464, 616),
    @(38, 13432, 962, 1324),
    @(39, 351, 581, 96),
    @(40, 574, 259, 230),
    @(41, 1834, 418, 68),
    @(42, 7228, 335, 263),
    @(43, 170, 145, 13),
    @(44, 514, 111, 49),
    @(45, 0, 0, 0),
    @(46, 1576, 601, 220),
    @(47, 7049, 1905, 1308),
    @(48, 16156, 1906, 2886),
    @(49, 7053, 1858, 874),
    @(50, 5939, 779, 1353),
    @(51, 16416, 1534, 2855),
    @(52, 2371, 534, 666),
    @(53, 5622, 1153, 1003),
    @(54, 1037, 713, 425),
    @(55, 1253, 732, 62),
    @(56, 2500, 665, 907),
    @(57, 6380, 2574, 1538),
    @(58, 9746, 1009, 393),
    @(59, 342533, 55531, 51621)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Restore the active cell/selection recorded in the workbook view
$ws.Range("B9").Select()
